$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Expansion")

# Rename the existing "Package Charges" display text to "Provider Package Charges"
$ws.Range("E10").Value = "Provider Package Charges"

# --- Row 11: HBP / NHA Package Charges ---
$ws.Range("A11").Formula = "=""1"""
$ws.Range("B11").Value = "http://hcxprotocol.io/codes/service-category"
$ws.Range("D11").Value = "HBP"
$ws.Range("E11").Value = "NHA Package Charges"
$ws.Range("F11").Formula = "=""false"""
$ws.Range("G11").Formula = "=""false"""

# --- Row 12: MJPJAY / Mahatma Jyotibha Phule PMJAY Package Charges ---
$ws.Range("A12").Formula = "=""1"""
$ws.Range("B12").Value = "http://hcxprotocol.io/codes/service-category"
$ws.Range("D12").Value = "MJPJAY"
$ws.Range("E12").Value = "Mahatma Jyotibha Phule PMJAY Package Charges"
$ws.Range("F12").Formula = "=""false"""
$ws.Range("G12").Formula = "=""false"""

# Collapse the helper formulas down to plain shared-string values
$ws.Range("A11:G12").Copy()
$ws.Range("A11:G12").PasteSpecial(-4163)

# Match the formatting used by the other data rows (row 2-10)
$ws.Range("A10:G10").Copy()
$ws.Range("A11:G12").PasteSpecial(-4122)
